$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected average experiment / agent-step timing values.
# For each data row, the previous "Avg_Agent_Step_Time" (G) and
# "Std_Agent_Step_Time" (M) values had mistakenly been left equal to the
# prior per-round values; they are now recalculated to their correct
# (larger) values, and "Avg_Experiment_Time" (H) / "Std_Experiment_Time"
# (N) are recalculated accordingly.

$ws.Range("G2").Value = 6.41424863
$ws.Range("H2").Value = 344.48750569
$ws.Range("M2").Value = 0.8515212965645228
$ws.Range("N2").Value = 88.38218234857277

$ws.Range("G3").Value = 10.97185791
$ws.Range("H3").Value = 984.2515611700002
$ws.Range("M3").Value = 1.549100725218475
$ws.Range("N3").Value = 275.9484917322532

$ws.Range("G4").Value = 1.86627926
$ws.Range("H4").Value = 55.36451846999999
$ws.Range("M4").Value = 0.3353193537182583
$ws.Range("N4").Value = 20.95915476158065

$ws.Range("G5").Value = 2.988098589999999
$ws.Range("H5").Value = 149.2233457
$ws.Range("M5").Value = 0.5331656409052948
$ws.Range("N5").Value = 54.40949202101353

$ws.Range("G6").Value = 0.49024171
$ws.Range("H6").Value = 8.056200349999999
$ws.Range("M6").Value = 0.1260675616210172
$ws.Range("N6").Value = 4.437641135113545

$ws.Range("G7").Value = 0.82102424
$ws.Range("H7").Value = 22.01964491
$ws.Range("M7").Value = 0.1790995858053783
$ws.Range("N7").Value = 10.44246761821167

$ws.Range("G8").Value = 0.21918958
$ws.Range("H8").Value = 2.33280048
$ws.Range("M8").Value = 0.05426064872468467
$ws.Range("N8").Value = 1.248818513056049

$ws.Range("G9").Value = 0.3740235
$ws.Range("H9").Value = 6.99994068
$ws.Range("M9").Value = 0.09163754644614111
$ws.Range("N9").Value = 3.785971300709705

$ws.Range("G10").Value = 0.11581932
$ws.Range("H10").Value = 0.9122120500000001
$ws.Range("M10").Value = 0.0270031369982331
$ws.Range("N10").Value = 0.440014622725552

$ws.Range("G11").Value = 0.20964824
$ws.Range("H11").Value = 3.08285206
$ws.Range("M11").Value = 0.0589490348301803
$ws.Range("N11").Value = 1.992381930183537

$ws.Range("G12").Value = 0.07555055999999999
$ws.Range("H12").Value = 0.48470289
$ws.Range("M12").Value = 0.01943454797368185
$ws.Range("N12").Value = 0.2535836909713923

$ws.Range("G13").Value = 0.13468841
$ws.Range("H13").Value = 1.59421345
$ws.Range("M13").Value = 0.03798739386776572
$ws.Range("N13").Value = 0.9918099429284268

$wb.Save()
